# "Excel for CompSci.xlsx" — fix a handful of mis-typed RGB / character
# values in the lookup table on Sheet1 (column B), and update the
# window's scroll/selection state to where the author was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B corrections -------------------------------------------------
# Several rows had the wrong "next" value typed in (duplicating the row
# above instead of incrementing by 10). Correct each one to the intended
# value.

$ws.Range("B8").Value  = "0,0,110"     # was a duplicate of B7 (0,0,100)
$ws.Range("B10").Value = "0,0,130"     # was a duplicate of B9 (0,0,120)
$ws.Range("B17").Value = "0,110,0"     # was a duplicate of B18 (0,120,0)
$ws.Range("B39").Value = 110110110     # was a duplicate of B40 (120120120)
$ws.Range("B40").Value = 120120120     # was a duplicate of B41... shift down
$ws.Range("B47").Value = "110,85,0"    # was 120,85,0
$ws.Range("B55").Value = "115,190,0"   # was 125,190,0

# --- Window state: selection / scroll position ----------------------------
# The author had scrolled down to row ~61 and left the selection on B56
# when the file was last saved (previously topLeftCell A63 / selection C89).
$ws.Range("B56").Select()
$excel.ActiveWindow.ScrollRow = 61
$excel.ActiveWindow.ScrollColumn = 1
